# Auto-generated Excel COM-interop script applying the Anima_Profits profit-table refresh
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 3000
$ws.Cells.Item(98, 9).Value = 0
$ws.Cells.Item(98, 11).Value = 0
$ws.Cells.Item(98, 13).ClearContents()
$ws.Cells.Item(107, 8).Value = 883.6129
$ws.Cells.Item(107, 9).Value = 1120.909
$ws.Cells.Item(107, 10).Value = 303.55554
$ws.Cells.Item(107, 11).Value = 1120.909
$ws.Cells.Item(107, 12).Value = 303.55554
$ws.Cells.Item(107, 13).Value = 799.0909999999999
$ws.Cells.Item(107, 14).Value = -4143.55554
$ws.Cells.Item(112, 8).Value = 5196.8613
$ws.Cells.Item(112, 10).Value = 5632.9395
$ws.Cells.Item(112, 12).Value = 16898.8185
$ws.Cells.Item(112, 14).Value = -19114.8185
$ws.Cells.Item(116, 8).Value = 7265.75
$ws.Cells.Item(116, 9).Value = 9035.666999999999
$ws.Cells.Item(116, 10).Value = 1956
$ws.Cells.Item(116, 11).Value = 9035.666999999999
$ws.Cells.Item(116, 12).Value = 1956
$ws.Cells.Item(116, 13).Value = -5593.666999999999
$ws.Cells.Item(116, 14).Value = -8840
$ws.Cells.Item(122, 8).Value = 3000
$ws.Cells.Item(122, 9).Value = 0
$ws.Cells.Item(122, 11).Value = 0
$ws.Cells.Item(122, 13).ClearContents()
$ws.Cells.Item(123, 8).Value = 31400
$ws.Cells.Item(123, 10).Value = 31400
$ws.Cells.Item(123, 12).Value = 31400
$ws.Cells.Item(123, 14).Value = -41200
$ws.Cells.Item(125, 8).Value = 1763.4
$ws.Cells.Item(125, 9).Value = 900
$ws.Cells.Item(125, 10).Value = 1859.3334
$ws.Cells.Item(125, 11).Value = 8100
$ws.Cells.Item(125, 12).Value = 16734.0006
$ws.Cells.Item(125, 13).Value = -5640
$ws.Cells.Item(125, 14).Value = -21654.0006
$ws.Cells.Item(132, 8).Value = 2738.535
$ws.Cells.Item(132, 9).Value = 2578.95
$ws.Cells.Item(132, 10).Value = 4866.3335
$ws.Cells.Item(132, 11).Value = 7736.849999999999
$ws.Cells.Item(132, 12).Value = 14599.0005
$ws.Cells.Item(132, 13).Value = -5206.849999999999
$ws.Cells.Item(132, 14).Value = -19659.0005
$ws.Cells.Item(135, 8).Value = 1912.76
$ws.Cells.Item(135, 9).Value = 815.4375
$ws.Cells.Item(135, 10).Value = 3863.5557
$ws.Cells.Item(135, 11).Value = 7338.9375
$ws.Cells.Item(135, 12).Value = 34772.0013
$ws.Cells.Item(135, 13).Value = -4803.9375
$ws.Cells.Item(135, 14).Value = -39842.0013
$ws.Cells.Item(138, 8).Value = 1405.58
$ws.Cells.Item(138, 9).Value = 638.6829
$ws.Cells.Item(138, 10).Value = 1938.5084
$ws.Cells.Item(138, 11).Value = 1916.0487
$ws.Cells.Item(138, 12).Value = 5815.5252
$ws.Cells.Item(138, 13).Value = 3223.9513
$ws.Cells.Item(138, 14).Value = -16095.5252

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(28, 8).Value = 13290.9
$ws.Cells.Item(28, 9).Value = 3821.111
$ws.Cells.Item(28, 11).Value = 3821.111
$ws.Cells.Item(28, 13).Value = -3629.111
$ws.Cells.Item(61, 8).Value = 2425.3215
$ws.Cells.Item(61, 9).Value = 2266.4583
$ws.Cells.Item(61, 11).Value = 2266.4583
$ws.Cells.Item(61, 13).Value = -2054.4583
$ws.Cells.Item(74, 8).Value = 1903.6666
$ws.Cells.Item(74, 9).Value = 1320.3636
$ws.Cells.Item(74, 10).Value = 2820.2856
$ws.Cells.Item(74, 11).Value = 1320.3636
$ws.Cells.Item(74, 12).Value = 2820.2856
$ws.Cells.Item(74, 13).Value = -446.3635999999999
$ws.Cells.Item(74, 14).Value = -4568.2856
$ws.Cells.Item(77, 8).Value = 1903.6666
$ws.Cells.Item(77, 9).Value = 1320.3636
$ws.Cells.Item(77, 10).Value = 2820.2856
$ws.Cells.Item(77, 11).Value = 6601.817999999999
$ws.Cells.Item(77, 12).Value = 14101.428
$ws.Cells.Item(77, 13).Value = -2233.817999999999
$ws.Cells.Item(77, 14).Value = -22837.428
$ws.Cells.Item(99, 8).Value = 13290.9
$ws.Cells.Item(99, 9).Value = 3821.111
$ws.Cells.Item(99, 11).Value = 3821.111
$ws.Cells.Item(99, 13).Value = -826.1109999999999
$ws.Cells.Item(102, 8).Value = 2231
$ws.Cells.Item(102, 9).Value = 1851.25
$ws.Cells.Item(102, 10).Value = 3750
$ws.Cells.Item(102, 11).Value = 1851.25
$ws.Cells.Item(102, 12).Value = 3750
$ws.Cells.Item(102, 13).Value = -229.25
$ws.Cells.Item(102, 14).Value = -6994
$ws.Cells.Item(132, 8).Value = 3070.1396
$ws.Cells.Item(132, 9).Value = 1809.871
$ws.Cells.Item(132, 10).Value = 6325.8335
$ws.Cells.Item(132, 11).Value = 5429.613
$ws.Cells.Item(132, 12).Value = 18977.5005
$ws.Cells.Item(132, 13).Value = -2899.613
$ws.Cells.Item(132, 14).Value = -24037.5005
$ws.Cells.Item(136, 8).Value = 2425.3215
$ws.Cells.Item(136, 9).Value = 2266.4583
$ws.Cells.Item(136, 11).Value = 6799.374899999999
$ws.Cells.Item(136, 13).Value = -4249.374899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 1356.5
$ws.Cells.Item(80, 9).Value = 2467.875
$ws.Cells.Item(80, 10).Value = 245.125
$ws.Cells.Item(80, 11).Value = 2467.875
$ws.Cells.Item(80, 12).Value = 245.125
$ws.Cells.Item(80, 13).Value = -1469.875
$ws.Cells.Item(80, 14).Value = -2241.125
$ws.Cells.Item(83, 8).Value = 1356.5
$ws.Cells.Item(83, 9).Value = 2467.875
$ws.Cells.Item(83, 10).Value = 245.125
$ws.Cells.Item(83, 11).Value = 12339.375
$ws.Cells.Item(83, 12).Value = 1225.625
$ws.Cells.Item(83, 13).Value = -7347.375
$ws.Cells.Item(83, 14).Value = -11209.625
$ws.Cells.Item(105, 8).Value = 41668770
$ws.Cells.Item(105, 9).Value = 41668770
$ws.Cells.Item(105, 10).Value = 0
$ws.Cells.Item(105, 11).Value = 41668770
$ws.Cells.Item(105, 12).Value = 0
$ws.Cells.Item(105, 13).Value = -41667023
$ws.Cells.Item(105, 14).ClearContents()
$ws.Cells.Item(134, 8).Value = 2531.3618
$ws.Cells.Item(134, 9).Value = 2117.8708
$ws.Cells.Item(134, 10).Value = 3332.5
$ws.Cells.Item(134, 11).Value = 6353.6124
$ws.Cells.Item(134, 12).Value = 9997.5
$ws.Cells.Item(134, 13).Value = -3818.6124
$ws.Cells.Item(134, 14).Value = -15067.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 1074
$ws.Cells.Item(58, 9).Value = 837.6786
$ws.Cells.Item(58, 10).Value = 1404.85
$ws.Cells.Item(58, 11).Value = 837.6786
$ws.Cells.Item(58, 12).Value = 1404.85
$ws.Cells.Item(58, 13).Value = -634.6786
$ws.Cells.Item(58, 14).Value = -1810.85
$ws.Cells.Item(99, 8).Value = 1995.238
$ws.Cells.Item(99, 9).Value = 1900
$ws.Cells.Item(99, 11).Value = 1900
$ws.Cells.Item(99, 13).Value = -402
$ws.Cells.Item(107, 8).Value = 929.13336
$ws.Cells.Item(107, 9).Value = 545.75
$ws.Cells.Item(107, 10).Value = 1068.5454
$ws.Cells.Item(107, 11).Value = 545.75
$ws.Cells.Item(107, 12).Value = 1068.5454
$ws.Cells.Item(107, 13).Value = 1374.25
$ws.Cells.Item(107, 14).Value = -4908.5454
$ws.Cells.Item(126, 8).Value = 1995.238
$ws.Cells.Item(126, 9).Value = 1900
$ws.Cells.Item(126, 11).Value = 5700
$ws.Cells.Item(126, 13).Value = -3230
$ws.Cells.Item(132, 8).Value = 1396.1163
$ws.Cells.Item(132, 9).Value = 1119.3422
$ws.Cells.Item(132, 11).Value = 3358.0266
$ws.Cells.Item(132, 13).Value = -828.0266000000001
$ws.Cells.Item(134, 8).Value = 4841.241
$ws.Cells.Item(134, 9).Value = 5414.4165
$ws.Cells.Item(134, 10).Value = 2090
$ws.Cells.Item(134, 11).Value = 16243.2495
$ws.Cells.Item(134, 12).Value = 6270
$ws.Cells.Item(134, 13).Value = -13708.2495
$ws.Cells.Item(134, 14).Value = -11340
$ws.Cells.Item(136, 8).Value = 1074
$ws.Cells.Item(136, 9).Value = 837.6786
$ws.Cells.Item(136, 10).Value = 1404.85
$ws.Cells.Item(136, 11).Value = 2513.0358
$ws.Cells.Item(136, 12).Value = 4214.549999999999
$ws.Cells.Item(136, 13).Value = 36.96420000000035
$ws.Cells.Item(136, 14).Value = -9314.549999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(97, 8).Value = 1239.8
$ws.Cells.Item(97, 9).Value = 839.6
$ws.Cells.Item(97, 10).Value = 1640
$ws.Cells.Item(97, 11).Value = 2518.8
$ws.Cells.Item(97, 12).Value = 4920
$ws.Cells.Item(97, 13).Value = -2022.8
$ws.Cells.Item(97, 14).Value = -5912
$ws.Cells.Item(125, 8).Value = 3792.7273
$ws.Cells.Item(125, 9).Value = 896.6667
$ws.Cells.Item(125, 10).Value = 4878.75
$ws.Cells.Item(125, 11).Value = 2690.0001
$ws.Cells.Item(125, 12).Value = 14636.25
$ws.Cells.Item(125, 13).Value = 2229.9999
$ws.Cells.Item(125, 14).Value = -24476.25
$ws.Cells.Item(131, 8).Value = 4107.5405
$ws.Cells.Item(131, 10).Value = 5868.4
$ws.Cells.Item(131, 12).Value = 17605.2
$ws.Cells.Item(131, 14).Value = -27685.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 1177378.8
$ws.Cells.Item(80, 9).Value = 1502671.6
$ws.Cells.Item(80, 10).Value = 201500
$ws.Cells.Item(80, 11).Value = 1502671.6
$ws.Cells.Item(80, 12).Value = 201500
$ws.Cells.Item(80, 13).Value = -1501673.6
$ws.Cells.Item(80, 14).Value = -203496
$ws.Cells.Item(83, 8).Value = 1177378.8
$ws.Cells.Item(83, 9).Value = 1502671.6
$ws.Cells.Item(83, 10).Value = 201500
$ws.Cells.Item(83, 11).Value = 7513358
$ws.Cells.Item(83, 12).Value = 1007500
$ws.Cells.Item(83, 13).Value = -7508366
$ws.Cells.Item(83, 14).Value = -1017484
$ws.Cells.Item(104, 8).Value = 33000
$ws.Cells.Item(104, 10).Value = 33000
$ws.Cells.Item(104, 12).Value = 33000
$ws.Cells.Item(104, 14).Value = -39988
$ws.Cells.Item(123, 8).Value = 10281.733
$ws.Cells.Item(123, 10).Value = 10281.733
$ws.Cells.Item(123, 12).Value = 10281.733
$ws.Cells.Item(123, 14).Value = -15181.733
$ws.Cells.Item(126, 8).Value = 3729.8572
$ws.Cells.Item(126, 9).Value = 3668.1667
$ws.Cells.Item(126, 10).Value = 4100
$ws.Cells.Item(126, 11).Value = 11004.5001
$ws.Cells.Item(126, 12).Value = 12300
$ws.Cells.Item(126, 13).Value = -8534.500100000001
$ws.Cells.Item(126, 14).Value = -17240
$ws.Cells.Item(132, 8).Value = 3613.25
$ws.Cells.Item(132, 9).Value = 3227.9355
$ws.Cells.Item(132, 10).Value = 6002.2
$ws.Cells.Item(132, 11).Value = 9683.806500000001
$ws.Cells.Item(132, 12).Value = 18006.6
$ws.Cells.Item(132, 13).Value = -7153.806500000001
$ws.Cells.Item(132, 14).Value = -23066.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 1608.1818
$ws.Cells.Item(68, 9).Value = 1698
$ws.Cells.Item(68, 10).Value = 1533.3334
$ws.Cells.Item(68, 11).Value = 1698
$ws.Cells.Item(68, 12).Value = 1533.3334
$ws.Cells.Item(68, 13).Value = -949
$ws.Cells.Item(68, 14).Value = -3031.3334
$ws.Cells.Item(71, 8).Value = 1608.1818
$ws.Cells.Item(71, 9).Value = 1698
$ws.Cells.Item(71, 10).Value = 1533.3334
$ws.Cells.Item(71, 11).Value = 8490
$ws.Cells.Item(71, 12).Value = 7666.666999999999
$ws.Cells.Item(71, 13).Value = -4746
$ws.Cells.Item(71, 14).Value = -15154.667
$ws.Cells.Item(122, 8).Value = 3562.0833
$ws.Cells.Item(122, 9).Value = 2916.6667
$ws.Cells.Item(122, 10).Value = 4207.5
$ws.Cells.Item(122, 11).Value = 8750.000100000001
$ws.Cells.Item(122, 12).Value = 12622.5
$ws.Cells.Item(122, 13).Value = -6300.000100000001
$ws.Cells.Item(122, 14).Value = -17522.5
$ws.Cells.Item(132, 8).Value = 3066.0732
$ws.Cells.Item(132, 9).Value = 2818.5557
$ws.Cells.Item(132, 10).Value = 3543.4285
$ws.Cells.Item(132, 11).Value = 8455.667099999999
$ws.Cells.Item(132, 12).Value = 10630.2855
$ws.Cells.Item(132, 13).Value = -5925.667099999999
$ws.Cells.Item(132, 14).Value = -15690.2855
$ws.Cells.Item(136, 8).Value = 5748982
$ws.Cells.Item(136, 9).Value = 1831.96
$ws.Cells.Item(136, 11).Value = 5495.88
$ws.Cells.Item(136, 13).Value = -2945.88

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(133, 8).Value = 40000
$ws.Cells.Item(133, 10).Value = 40000
$ws.Cells.Item(133, 12).Value = 40000
$ws.Cells.Item(133, 14).Value = -50120
$ws.Cells.Item(136, 8).Value = 2062.096
$ws.Cells.Item(136, 9).Value = 1741.5862
$ws.Cells.Item(136, 10).Value = 3301.4
$ws.Cells.Item(136, 11).Value = 5224.7586
$ws.Cells.Item(136, 12).Value = 9904.200000000001
$ws.Cells.Item(136, 13).Value = -2674.7586
$ws.Cells.Item(136, 14).Value = -15004.2
